# Auto-generated Excel COM-interop script to apply the scheduled-runner Leve profit updates
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H18").Value = 697.6667
$ws_ALC.Range("I18").Value = 697.2
$ws_ALC.Range("J18").Value = 700
$ws_ALC.Range("K18").Value = 697.2
$ws_ALC.Range("L18").Value = 700
$ws_ALC.Range("M18").Value = -413.2
$ws_ALC.Range("N18").Value = -1268
$ws_ALC.Range("H40").Value = 1501.4286
$ws_ALC.Range("I40").Value = 1400
$ws_ALC.Range("J40").Value = 1549.4736
$ws_ALC.Range("K40").Value = 1400
$ws_ALC.Range("L40").Value = 1549.4736
$ws_ALC.Range("M40").Value = -1225
$ws_ALC.Range("N40").Value = -1899.4736
$ws_ALC.Range("H64").Value = 2791.9355
$ws_ALC.Range("I64").Value = 2694.2856
$ws_ALC.Range("J64").Value = 2918.5186
$ws_ALC.Range("K64").Value = 2694.2856
$ws_ALC.Range("L64").Value = 2918.5186
$ws_ALC.Range("M64").Value = -2446.2856
$ws_ALC.Range("N64").Value = -3414.5186
$ws_ALC.Range("H67").Value = 2791.9355
$ws_ALC.Range("I67").Value = 2694.2856
$ws_ALC.Range("J67").Value = 2918.5186
$ws_ALC.Range("K67").Value = 2694.2856
$ws_ALC.Range("L67").Value = 2918.5186
$ws_ALC.Range("M67").Value = -1836.2856
$ws_ALC.Range("N67").Value = -4634.518599999999
$ws_ALC.Range("H131").Value = 1370.2142
$ws_ALC.Range("I131").Value = 1081.9166
$ws_ALC.Range("J131").Value = 3100
$ws_ALC.Range("K131").Value = 3245.7498
$ws_ALC.Range("L131").Value = 9300
$ws_ALC.Range("M131").Value = 1794.2502
$ws_ALC.Range("N131").Value = -19380
$ws_ALC.Range("H137").Value = 2562.7273
$ws_ALC.Range("I137").Value = 1965
$ws_ALC.Range("J137").Value = 3280
$ws_ALC.Range("K137").Value = 5895
$ws_ALC.Range("L137").Value = 9840
$ws_ALC.Range("M137").Value = -3345
$ws_ALC.Range("N137").Value = -14940
$ws_ALC.Range("H138").Value = 2278.9124
$ws_ALC.Range("I138").Value = 2365.4285
$ws_ALC.Range("J138").Value = 2250.7441
$ws_ALC.Range("K138").Value = 7096.2855
$ws_ALC.Range("L138").Value = 6752.2323
$ws_ALC.Range("M138").Value = -1956.2855
$ws_ALC.Range("N138").Value = -17032.2323

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 1443.7142
$ws_ARM.Range("I61").Value = 1355.2307
$ws_ARM.Range("J61").Value = 1587.5
$ws_ARM.Range("K61").Value = 1355.2307
$ws_ARM.Range("L61").Value = 1587.5
$ws_ARM.Range("M61").Value = -1143.2307
$ws_ARM.Range("N61").Value = -2011.5
$ws_ARM.Range("H74").Value = 1064.7742
$ws_ARM.Range("I74").Value = 779.04346
$ws_ARM.Range("J74").Value = 1886.25
$ws_ARM.Range("K74").Value = 779.04346
$ws_ARM.Range("L74").Value = 1886.25
$ws_ARM.Range("M74").Value = 94.95654000000002
$ws_ARM.Range("N74").Value = -3634.25
$ws_ARM.Range("H77").Value = 1064.7742
$ws_ARM.Range("I77").Value = 779.04346
$ws_ARM.Range("J77").Value = 1886.25
$ws_ARM.Range("K77").Value = 3895.2173
$ws_ARM.Range("L77").Value = 9431.25
$ws_ARM.Range("M77").Value = 472.7827000000002
$ws_ARM.Range("N77").Value = -18167.25
$ws_ARM.Range("H136").Value = 1443.7142
$ws_ARM.Range("I136").Value = 1355.2307
$ws_ARM.Range("J136").Value = 1587.5
$ws_ARM.Range("K136").Value = 4065.6921
$ws_ARM.Range("L136").Value = 4762.5
$ws_ARM.Range("M136").Value = -1515.6921
$ws_ARM.Range("N136").Value = -9862.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H22").Value = 554
$ws_BSM.Range("I22").Value = 600
$ws_BSM.Range("J22").Value = 492.66666
$ws_BSM.Range("K22").Value = 600
$ws_BSM.Range("L22").Value = 492.66666
$ws_BSM.Range("M22").Value = -427
$ws_BSM.Range("N22").Value = -838.66666
$ws_BSM.Range("H122").Value = 29408.46
$ws_BSM.Range("J122").Value = 29408.46
$ws_BSM.Range("L122").Value = 29408.46
$ws_BSM.Range("N122").Value = -39208.46
$ws_BSM.Range("H124").Value = 44640
$ws_BSM.Range("J124").Value = 44640
$ws_BSM.Range("L124").Value = 44640
$ws_BSM.Range("N124").Value = -54460
$ws_BSM.Range("H134").Value = 1256.5385
$ws_BSM.Range("I134").Value = 1250.3636
$ws_BSM.Range("J134").Value = 1290.5
$ws_BSM.Range("K134").Value = 3751.0908
$ws_BSM.Range("L134").Value = 3871.5
$ws_BSM.Range("M134").Value = -1216.0908
$ws_BSM.Range("N134").Value = -8941.5

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 2636.697
$ws_CRP.Range("I31").Value = 1795.15
$ws_CRP.Range("J31").Value = 3931.3845
$ws_CRP.Range("K31").Value = 1795.15
$ws_CRP.Range("L31").Value = 3931.3845
$ws_CRP.Range("M31").Value = -1500.15
$ws_CRP.Range("N31").Value = -4521.3845
$ws_CRP.Range("H34").Value = 2636.697
$ws_CRP.Range("I34").Value = 1795.15
$ws_CRP.Range("J34").Value = 3931.3845
$ws_CRP.Range("K34").Value = 1795.15
$ws_CRP.Range("L34").Value = 3931.3845
$ws_CRP.Range("M34").Value = -1593.15
$ws_CRP.Range("N34").Value = -4335.3845
$ws_CRP.Range("H42").Value = 3997.5
$ws_CRP.Range("J42").Value = 3997.5
$ws_CRP.Range("L42").Value = 3997.5
$ws_CRP.Range("N42").Value = -5183.5
$ws_CRP.Range("H55").Value = 5208.6665
$ws_CRP.Range("I55").Value = 4364.5
$ws_CRP.Range("J55").Value = 5630.75
$ws_CRP.Range("K55").Value = 4364.5
$ws_CRP.Range("L55").Value = 5630.75
$ws_CRP.Range("M55").Value = -4049.5
$ws_CRP.Range("N55").Value = -6260.75
$ws_CRP.Range("H58").Value = 1379.7778
$ws_CRP.Range("I58").Value = 1326.4572
$ws_CRP.Range("J58").Value = 1566.4
$ws_CRP.Range("K58").Value = 1326.4572
$ws_CRP.Range("L58").Value = 1566.4
$ws_CRP.Range("M58").Value = -1123.4572
$ws_CRP.Range("N58").Value = -1972.4
$ws_CRP.Range("H62").Value = 3599
$ws_CRP.Range("I62").Value = 3538.75
$ws_CRP.Range("J62").Value = 3840
$ws_CRP.Range("K62").Value = 3538.75
$ws_CRP.Range("L62").Value = 3840
$ws_CRP.Range("M62").Value = -2914.75
$ws_CRP.Range("N62").Value = -5088
$ws_CRP.Range("H65").Value = 3599
$ws_CRP.Range("I65").Value = 3538.75
$ws_CRP.Range("J65").Value = 3840
$ws_CRP.Range("K65").Value = 17693.75
$ws_CRP.Range("L65").Value = 19200
$ws_CRP.Range("M65").Value = -14573.75
$ws_CRP.Range("N65").Value = -25440
$ws_CRP.Range("H136").Value = 1379.7778
$ws_CRP.Range("I136").Value = 1326.4572
$ws_CRP.Range("J136").Value = 1566.4
$ws_CRP.Range("K136").Value = 3979.3716
$ws_CRP.Range("L136").Value = 4699.200000000001
$ws_CRP.Range("M136").Value = -1429.3716
$ws_CRP.Range("N136").Value = -9799.200000000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H80").Value = 8175
$ws_CUL.Range("I80").Value = 0
$ws_CUL.Range("J80").Value = 8175
$ws_CUL.Range("K80").Value = 0
$ws_CUL.Range("L80").Value = 24525
$ws_CUL.Range("M80").ClearContents()
$ws_CUL.Range("N80").Value = -26397
$ws_CUL.Range("H83").Value = 8175
$ws_CUL.Range("I83").Value = 0
$ws_CUL.Range("J83").Value = 8175
$ws_CUL.Range("K83").Value = 0
$ws_CUL.Range("L83").Value = 73575
$ws_CUL.Range("M83").ClearContents()
$ws_CUL.Range("N83").Value = -82935

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H54").Value = 7725
$ws_GSM.Range("J54").Value = 7725
$ws_GSM.Range("L54").Value = 7725
$ws_GSM.Range("N54").Value = -8505
$ws_GSM.Range("H58").Value = 0
$ws_GSM.Range("I58").Value = 0
$ws_GSM.Range("K58").Value = 0
$ws_GSM.Range("M58").ClearContents()
$ws_GSM.Range("H124").Value = 60000
$ws_GSM.Range("J124").Value = 60000
$ws_GSM.Range("L124").Value = 60000
$ws_GSM.Range("N124").Value = -69820
$ws_GSM.Range("H126").Value = 10420.387
$ws_GSM.Range("I126").Value = 2247.1765
$ws_GSM.Range("J126").Value = 20345
$ws_GSM.Range("K126").Value = 6741.529500000001
$ws_GSM.Range("L126").Value = 61035
$ws_GSM.Range("M126").Value = -4271.529500000001
$ws_GSM.Range("N126").Value = -65975
$ws_GSM.Range("H132").Value = 2997.8696
$ws_GSM.Range("I132").Value = 2249.4443
$ws_GSM.Range("J132").Value = 3479
$ws_GSM.Range("K132").Value = 6748.3329
$ws_GSM.Range("L132").Value = 10437
$ws_GSM.Range("M132").Value = -4218.3329
$ws_GSM.Range("N132").Value = -15497

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H68").Value = 1618
$ws_LTW.Range("I68").Value = 1618
$ws_LTW.Range("K68").Value = 1618
$ws_LTW.Range("M68").Value = -869
$ws_LTW.Range("H71").Value = 1618
$ws_LTW.Range("I71").Value = 1618
$ws_LTW.Range("K71").Value = 8090
$ws_LTW.Range("M71").Value = -4346
$ws_LTW.Range("H100").Value = 1060.3846
$ws_LTW.Range("I100").Value = 1060.3846
$ws_LTW.Range("K100").Value = 1060.3846
$ws_LTW.Range("M100").Value = -519.3846000000001
$ws_LTW.Range("H122").Value = 3159.5
$ws_LTW.Range("I122").Value = 2373.875
$ws_LTW.Range("J122").Value = 3945.125
$ws_LTW.Range("K122").Value = 7121.625
$ws_LTW.Range("L122").Value = 11835.375
$ws_LTW.Range("M122").Value = -4671.625
$ws_LTW.Range("N122").Value = -16735.375
$ws_LTW.Range("H132").Value = 3315.639
$ws_LTW.Range("I132").Value = 3651.5
$ws_LTW.Range("J132").Value = 2895.8125
$ws_LTW.Range("K132").Value = 10954.5
$ws_LTW.Range("L132").Value = 8687.4375
$ws_LTW.Range("M132").Value = -8424.5
$ws_LTW.Range("N132").Value = -13747.4375
$ws_LTW.Range("H136").Value = 2192.25
$ws_LTW.Range("I136").Value = 1715.2609
$ws_LTW.Range("J136").Value = 3036.1538
$ws_LTW.Range("K136").Value = 5145.7827
$ws_LTW.Range("L136").Value = 9108.4614
$ws_LTW.Range("M136").Value = -2595.7827
$ws_LTW.Range("N136").Value = -14208.4614

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H136").Value = 908.55817
$ws_WVR.Range("I136").Value = 723.84375
$ws_WVR.Range("J136").Value = 1445.909
$ws_WVR.Range("K136").Value = 2171.53125
$ws_WVR.Range("L136").Value = 4337.727000000001
$ws_WVR.Range("M136").Value = 378.46875
$ws_WVR.Range("N136").Value = -9437.727000000001
